$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 811.1429000000001
$ws.Range("I4").Value = 696.3333
$ws.Range("K4").Value = 696.3333
$ws.Range("M4").Value = -582.3333
$ws.Range("H40").Value = 2883.7827
$ws.Range("J40").Value = 2872.8823
$ws.Range("L40").Value = 2872.8823
$ws.Range("N40").Value = -3222.8823
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H61").Value = 6665.3335
$ws.Range("I61").Value = 6665.3335
$ws.Range("K61").Value = 19996.0005
$ws.Range("M61").Value = -19824.0005
$ws.Range("H62").Value = 2414.6667
$ws.Range("I62").Value = 2438.7058
$ws.Range("K62").Value = 2438.7058
$ws.Range("M62").Value = -1814.7058
$ws.Range("H65").Value = 2414.6667
$ws.Range("I65").Value = 2438.7058
$ws.Range("K65").Value = 12193.529
$ws.Range("M65").Value = -9073.529
$ws.Range("H86").Value = 6057.278
$ws.Range("I86").Value = 5723.5713
$ws.Range("J86").Value = 7225.25
$ws.Range("K86").Value = 5723.5713
$ws.Range("L86").Value = 7225.25
$ws.Range("M86").Value = -4600.5713
$ws.Range("N86").Value = -9471.25
$ws.Range("H89").Value = 6057.278
$ws.Range("I89").Value = 5723.5713
$ws.Range("J89").Value = 7225.25
$ws.Range("K89").Value = 28617.8565
$ws.Range("L89").Value = 36126.25
$ws.Range("M89").Value = -23001.8565
$ws.Range("N89").Value = -47358.25
$ws.Range("H92").Value = 931.5
$ws.Range("I92").Value = 958.5925999999999
$ws.Range("J92").Value = 200
$ws.Range("K92").Value = 958.5925999999999
$ws.Range("L92").Value = 200
$ws.Range("M92").Value = 289.4074000000001
$ws.Range("N92").Value = -2696
$ws.Range("H129").Value = 2387.5557
$ws.Range("I129").Value = 1581.3334
$ws.Range("K129").Value = 4744.0002
$ws.Range("M129").Value = 255.9997999999996
$ws.Range("H132").Value = 4051.5107
$ws.Range("I132").Value = 3600.5908
$ws.Range("K132").Value = 10801.7724
$ws.Range("M132").Value = -8271.7724
$ws.Range("H137").Value = 2703.1667
$ws.Range("I137").Value = 3014.1428
$ws.Range("K137").Value = 9042.428400000001
$ws.Range("M137").Value = -6492.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 40007204
$ws.Range("I74").Value = 4574.6924
$ws.Range("J74").Value = 83343384
$ws.Range("K74").Value = 4574.6924
$ws.Range("L74").Value = 83343384
$ws.Range("M74").Value = -3700.6924
$ws.Range("N74").Value = -83345132
$ws.Range("H77").Value = 40007204
$ws.Range("I77").Value = 4574.6924
$ws.Range("J77").Value = 83343384
$ws.Range("K77").Value = 22873.462
$ws.Range("L77").Value = 416716920
$ws.Range("M77").Value = -18505.462
$ws.Range("N77").Value = -416725656
$ws.Range("H110").Value = 1742
$ws.Range("I110").Value = 1843.5555
$ws.Range("J110").Value = 1132.6666
$ws.Range("K110").Value = 1843.5555
$ws.Range("L110").Value = 1132.6666
$ws.Range("M110").Value = 201.4445000000001
$ws.Range("N110").Value = -5222.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 24513584
$ws.Range("I20").Value = 30868588
$ws.Range("K20").Value = 30868588
$ws.Range("M20").Value = -30868341
$ws.Range("H86").Value = 2420.5527
$ws.Range("I86").Value = 2006.4333
$ws.Range("K86").Value = 2006.4333
$ws.Range("M86").Value = -883.4332999999999
$ws.Range("H89").Value = 2420.5527
$ws.Range("I89").Value = 2006.4333
$ws.Range("K89").Value = 10032.1665
$ws.Range("M89").Value = -4416.166499999999
$ws.Range("H105").Value = 11819709
$ws.Range("I105").Value = 626008.1
$ws.Range("J105").Value = 41669580
$ws.Range("K105").Value = 626008.1
$ws.Range("L105").Value = 41669580
$ws.Range("M105").Value = -624261.1
$ws.Range("N105").Value = -41673074
$ws.Range("H107").Value = 2748715.2
$ws.Range("I107").Value = 3078260
$ws.Range("K107").Value = 3078260
$ws.Range("M107").Value = -3076340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2909.0833
$ws.Range("I58").Value = 2363.75
$ws.Range("J58").Value = 3999.75
$ws.Range("K58").Value = 2363.75
$ws.Range("L58").Value = 3999.75
$ws.Range("M58").Value = -2160.75
$ws.Range("N58").Value = -4405.75
$ws.Range("H136").Value = 2909.0833
$ws.Range("I136").Value = 2363.75
$ws.Range("J136").Value = 3999.75
$ws.Range("K136").Value = 7091.25
$ws.Range("L136").Value = 11999.25
$ws.Range("M136").Value = -4541.25
$ws.Range("N136").Value = -17099.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 203.72728
$ws.Range("I2").Value = 137.92308
$ws.Range("J2").Value = 298.77777
$ws.Range("K2").Value = 827.5384799999999
$ws.Range("L2").Value = 1792.66662
$ws.Range("M2").Value = -714.5384799999999
$ws.Range("N2").Value = -2018.66662
$ws.Range("H38").Value = 575.1429000000001
$ws.Range("I38").Value = 153.5
$ws.Range("K38").Value = 460.5
$ws.Range("M38").Value = -113.5
$ws.Range("H121").Value = 5327088
$ws.Range("I121").Value = 16683956
$ws.Range("J121").Value = 85456.69500000001
$ws.Range("K121").Value = 50051868
$ws.Range("L121").Value = 256370.085
$ws.Range("M121").Value = -50050558
$ws.Range("N121").Value = -258990.085
$ws.Range("H129").Value = 80392.57000000001
$ws.Range("I129").Value = 937.5
$ws.Range("J129").Value = 186332.67
$ws.Range("K129").Value = 2812.5
$ws.Range("L129").Value = 558998.01
$ws.Range("M129").Value = 2187.5
$ws.Range("N129").Value = -568998.01
$ws.Range("H132").Value = 6816.6787
$ws.Range("I132").Value = 5941.222
$ws.Range("K132").Value = 53470.998
$ws.Range("M132").Value = -50940.998
$ws.Range("H133").Value = 3537
$ws.Range("I133").Value = 3537
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 10611
$ws.Range("L133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -5551
$ws.Range("H137").Value = 2896.4443
$ws.Range("I137").Value = 2909.125
$ws.Range("J137").Value = 2795
$ws.Range("K137").Value = 8727.375
$ws.Range("L137").Value = 8385
$ws.Range("M137").Value = -3627.375
$ws.Range("N137").Value = -18585

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 875000.3
$ws.Range("I14").Value = 875000.3
$ws.Range("K14").Value = 875000.3
$ws.Range("M14").Value = -874832.3
$ws.Range("H80").Value = 50003636
$ws.Range("J80").Value = 6125.25
$ws.Range("L80").Value = 6125.25
$ws.Range("N80").Value = -8121.25
$ws.Range("H83").Value = 50003636
$ws.Range("J83").Value = 6125.25
$ws.Range("L83").Value = 30626.25
$ws.Range("N83").Value = -40610.25
$ws.Range("H122").Value = 66672012
$ws.Range("I122").Value = 55559788
$ws.Range("K122").Value = 166679364
$ws.Range("M122").Value = -166676914
$ws.Range("H132").Value = 1596.6086
$ws.Range("I132").Value = 1466.8667
$ws.Range("J132").Value = 1839.875
$ws.Range("K132").Value = 4400.6001
$ws.Range("L132").Value = 5519.625
$ws.Range("M132").Value = -1870.6001
$ws.Range("N132").Value = -10579.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3323
$ws.Range("I22").Value = 1520
$ws.Range("K22").Value = 1520
$ws.Range("M22").Value = -1225
$ws.Range("H27").Value = 3323
$ws.Range("I27").Value = 1520
$ws.Range("K27").Value = 1520
$ws.Range("M27").Value = -1413
$ws.Range("H61").Value = 7200.1055
$ws.Range("I61").Value = 6635.4116
$ws.Range("J61").Value = 12000
$ws.Range("K61").Value = 6635.4116
$ws.Range("L61").Value = 12000
$ws.Range("M61").Value = -6433.4116
$ws.Range("N61").Value = -12404
$ws.Range("H100").Value = 2104.4
$ws.Range("I100").Value = 2309
$ws.Range("J100").Value = 1797.5
$ws.Range("K100").Value = 2309
$ws.Range("L100").Value = 1797.5
$ws.Range("M100").Value = -1768
$ws.Range("N100").Value = -2879.5
$ws.Range("H113").Value = 7200.1055
$ws.Range("I113").Value = 6635.4116
$ws.Range("J113").Value = 12000
$ws.Range("K113").Value = 6635.4116
$ws.Range("L113").Value = 12000
$ws.Range("M113").Value = -4465.4116
$ws.Range("N113").Value = -16340
$ws.Range("H122").Value = 5039.909
$ws.Range("I122").Value = 3413.1667
$ws.Range("J122").Value = 6992
$ws.Range("K122").Value = 10239.5001
$ws.Range("L122").Value = 20976
$ws.Range("M122").Value = -7789.500100000001
$ws.Range("N122").Value = -25876
$ws.Range("H132").Value = 7385.1396
$ws.Range("I132").Value = 7278.517
$ws.Range("J132").Value = 7606
$ws.Range("K132").Value = 21835.551
$ws.Range("L132").Value = 22818
$ws.Range("M132").Value = -19305.551
$ws.Range("N132").Value = -27878
$ws.Range("H136").Value = 6438.1304
$ws.Range("I136").Value = 4821
$ws.Range("K136").Value = 14463
$ws.Range("M136").Value = -11913

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 7856.5713
$ws.Range("I14").Value = 5000
$ws.Range("K14").Value = 5000
$ws.Range("M14").Value = -4832
$ws.Range("H107").Value = 531.5625
$ws.Range("I107").Value = 537.3333
$ws.Range("K107").Value = 1611.9999
$ws.Range("M107").Value = 308.0001
$ws.Range("H122").Value = 7355484.5
$ws.Range("I122").Value = 2580.4814
$ws.Range("J122").Value = 35716690
$ws.Range("K122").Value = 7741.4442
$ws.Range("L122").Value = 107150070
$ws.Range("M122").Value = -5291.4442
$ws.Range("N122").Value = -107154970
$ws.Range("H132").Value = 3669.75
$ws.Range("I132").Value = 3562.75
$ws.Range("K132").Value = 10688.25
$ws.Range("M132").Value = -8158.25
$ws.Range("H136").Value = 19618232
$ws.Range("I136").Value = 23819278
$ws.Range("K136").Value = 71457834
$ws.Range("M136").Value = -71455284
